$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) updates
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 47
$ws1.Range("F3").Value = 21356
$ws1.Range("F6").Value = 1133
$ws1.Range("F8").Value = 7956
$ws1.Range("F11").Value = 765
$ws1.Range("F12").Value = 317
$ws1.Range("F14").Value = 190
$ws1.Range("F15").Value = 174
$ws1.Range("F20").Value = 542
$ws1.Range("F27").Value = 1192
$ws1.Range("F28").Value = 56
$ws1.Range("F32").Value = 608
$ws1.Range("F34").Value = 143
$ws1.Range("F35").Value = 5076
$ws1.Range("F40").Value = 13160
$ws1.Range("F45").Value = 312

# Sheet "全部类型" (sheet4) updates
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 47
$ws4.Range("F3").Value = 21356
$ws4.Range("F5").Value = 1133
$ws4.Range("F7").Value = 7956
$ws4.Range("F10").Value = 765
$ws4.Range("F11").Value = 317
$ws4.Range("F13").Value = 190
$ws4.Range("F14").Value = 174
$ws4.Range("F18").Value = 542
$ws4.Range("F25").Value = 1192
$ws4.Range("F26").Value = 56
$ws4.Range("F30").Value = 608
$ws4.Range("F33").Value = 143
$ws4.Range("F35").Value = 5076
$ws4.Range("F40").Value = 13160
$ws4.Range("F45").Value = 312

